{"js": "// \"Version 2.\" -> \"Version 1.\" (wireframes.docx version bump / revert)\n//\n// The paragraph starts out as four runs: \"Versi\" | \"on\" | \" 2\" | \".\"\n// (with a spellStart/spellEnd proofErr pair wrapping \"Versi\"+\"on\").\n// The target keeps \"Version\" as a single run inside the proofErr pair,\n// replaces the \" 2\" run's text with \" 1.\", and drops the trailing \".\"\n// run entirely. We do this with three in-place text operations instead\n// of a single whole-paragraph rewrite so the existing run/proofErr\n// structure survives the edit.\n\nconst body = context.document.body;\n\n// Step 1: normalize \"Versi\" + \"on\" into a single \"Version\" run by\n// replacing the found (cross-run) match with the same text; Word's\n// insertText(\"Replace\") collapses the match into one run.\nconst versionMatches = body.search(\"Version\", { matchCase: true });\nversionMatches.load(\"items\");\nawait context.sync();\n\nif (versionMatches.items.length > 0) {\n  versionMatches.items[0].insertText(\"Version\", \"Replace\");\n  await context.sync();\n}\n\n// Step 2: change the \" 2\" run's text to \" 1.\" in place.\nconst numberMatches = context.document.body.search(\" 2\", { matchCase: true });\nnumberMatches.load(\"items\");\nawait context.sync();\n\nif (numberMatches.items.length > 0) {\n  numberMatches.items[0].insertText(\" 1.\", \"Replace\");\n  await context.sync();\n}\n\n// Step 3: remove the now-redundant trailing \".\" run (the original final\n// sentence period, now duplicated by the \" 1.\" run from Step 2).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\nconst dotMatches = paragraph.search(\".\", { matchCase: true });\ndotMatches.load(\"items\");\nawait context.sync();\n\nif (dotMatches.items.length > 0) {\n  dotMatches.items[dotMatches.items.length - 1].delete();\n  await context.sync();\n}\n", "ps1": "# \"Version 2.\" -> \"Version 1.\" (wireframes.docx version bump / revert)\n#\n# The paragraph starts out as four runs: \"Versi\" | \"on\" | \" 2\" | \".\"\n# (with a spellStart/spellEnd proofErr pair wrapping \"Versi\"+\"on\").\n# The target keeps \"Version\" as a single run inside the proofErr pair,\n# replaces the \" 2\" run's text with \" 1.\", and drops the trailing \".\"\n# run entirely. We do this with targeted Find/Replace + a single\n# character delete instead of rewriting the whole paragraph, so the\n# existing run/proofErr structure survives the edit.\n\n$d = $word.ActiveDocument\n\n# Step 1: normalize \"Versi\" + \"on\" into a single \"Version\" run by\n# running a find/replace over the (cross-run) match; Word collapses\n# the replaced text into one run.\n$d.Content.Find.Execute(\"Version\", $false, $false, $false, $false, $false, $true, 1, $false, \"Version\", 2) | Out-Null\n\n# Step 2: change the \" 2\" run's text to \" 1.\" in place.\n$d.Content.Find.Execute(\" 2\", $true, $false, $false, $false, $false, $true, 1, $false, \" 1.\", 2) | Out-Null\n\n# Step 3: remove the now-redundant trailing \".\" run (the original final\n# sentence period, now duplicated by the \" 1.\" text from Step 2). It is\n# the last text character of the paragraph, just before the paragraph\n# mark.\n$p = $d.Paragraphs(1)\n$pr = $p.Range\n$lastCharIndex = $pr.Characters.Count - 1\n$pr.Characters($lastCharIndex).Delete()\n"}
